$wb = $excel.ActiveWorkbook

# --- Rename header labels on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "PO Forecast"

# --- Header row ---
$ws.Range("A1").Value = "ds"
$ws.Range("B1").Value = "PO_Forecast"
$ws.Range("C1").Value = "yhat_lower"
$ws.Range("D1").Value = "yhat_upper"

# --- Data rows (2-23) ---
$ws.Range("A2").Value = 45473.99999999999
$ws.Range("B2").Value = 137
$ws.Range("C2").Value = -304.4307845335903
$ws.Range("D2").Value = 546.3838193236775
$ws.Range("A3").Value = 45487.99999999999
$ws.Range("B3").Value = 170
$ws.Range("C3").Value = -279.3898657164173
$ws.Range("D3").Value = 611.2520013179035
$ws.Range("A4").Value = 45494.99999999999
$ws.Range("B4").Value = 186
$ws.Range("C4").Value = -261.6357501063221
$ws.Range("D4").Value = 633.4213625978099
$ws.Range("A5").Value = 45508.99999999999
$ws.Range("B5").Value = 218
$ws.Range("C5").Value = -199.8640034953189
$ws.Range("D5").Value = 637.8107158910641
$ws.Range("A6").Value = 45515.99999999999
$ws.Range("B6").Value = 234
$ws.Range("C6").Value = -211.6062353474807
$ws.Range("D6").Value = 657.3909917609443
$ws.Range("A7").Value = 45522.99999999999
$ws.Range("B7").Value = 250
$ws.Range("C7").Value = -162.5952574308227
$ws.Range("D7").Value = 709.2435630845264
$ws.Range("A8").Value = 45536.99999999999
$ws.Range("B8").Value = 283
$ws.Range("C8").Value = -121.8465023084227
$ws.Range("D8").Value = 700.1009759535943
$ws.Range("A9").Value = 45543.99999999999
$ws.Range("B9").Value = 299
$ws.Range("C9").Value = -146.1659260207618
$ws.Range("D9").Value = 734.3242818162198
$ws.Range("A10").Value = 45550.99999999999
$ws.Range("B10").Value = 315
$ws.Range("C10").Value = -101.006277581766
$ws.Range("D10").Value = 738.1127763991782
$ws.Range("A11").Value = 45557.99999999999
$ws.Range("B11").Value = 331
$ws.Range("C11").Value = -100.1279509449282
$ws.Range("D11").Value = 747.2787739822336
$ws.Range("A12").Value = 45585.99999999999
$ws.Range("B12").Value = 396
$ws.Range("C12").Value = 16.71531126755297
$ws.Range("D12").Value = 799.7390785247901
$ws.Range("A13").Value = 45592.99999999999
$ws.Range("B13").Value = 412
$ws.Range("C13").Value = -15.38840843791121
$ws.Range("D13").Value = 826.6364122101511
$ws.Range("A14").Value = 45599.99999999999
$ws.Range("B14").Value = 428
$ws.Range("C14").Value = -34.74314219459435
$ws.Range("D14").Value = 833.9740334359805
$ws.Range("A15").Value = 45613.99999999999
$ws.Range("B15").Value = 461
$ws.Range("C15").Value = 49.49710866673044
$ws.Range("D15").Value = 895.1178116678025
$ws.Range("A16").Value = 45620.99999999999
$ws.Range("B16").Value = 477
$ws.Range("C16").Value = 74.73542065684317
$ws.Range("D16").Value = 897.261102573701
$ws.Range("A17").Value = 45627.99999999999
$ws.Range("B17").Value = 493
$ws.Range("C17").Value = 74.72781724546317
$ws.Range("D17").Value = 909.7433306217238
$ws.Range("A18").Value = 45634.99999999999
$ws.Range("B18").Value = 509
$ws.Range("C18").Value = 116.6429186475728
$ws.Range("D18").Value = 925.2884518720851
$ws.Range("A19").Value = 45641.99999999999
$ws.Range("B19").Value = 525
$ws.Range("C19").Value = 97.34844136047836
$ws.Range("D19").Value = 975.0285578880989
$ws.Range("A20").Value = 45648.99999999999
$ws.Range("B20").Value = 541
$ws.Range("C20").Value = 104.3663646870781
$ws.Range("D20").Value = 976.5229486052755
$ws.Range("A21").Value = 45655.99999999999
$ws.Range("B21").Value = 558
$ws.Range("C21").Value = 99.25521878760424
$ws.Range("D21").Value = 971.5211144668652
$ws.Range("A22").Value = 45662.99999999999
$ws.Range("B22").Value = 574
$ws.Range("C22").Value = 174.8454869857925
$ws.Range("D22").Value = 990.6500173269834
$ws.Range("A23").Value = 45669.99999999999
$ws.Range("B23").Value = 590
$ws.Range("C23").Value = 136.2058086749927
$ws.Range("D23").Value = 1026.30139014597

# --- Formatting: copy header style from existing sheet header cell to row 1 ---
$wsWeekly.Range("A1").Copy() | Out-Null
$ws.Range("A1:D1").PasteSpecial(-4122) | Out-Null

# --- Formatting: copy date style from existing sheet date cell to column A (rows 2-23) ---
$wsWeekly.Range("A2").Copy() | Out-Null
$ws.Range("A2:A23").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Restore original active sheet/selection so the workbook view is unaffected ---
$wsWeekly.Activate()
$wsWeekly.Range("A1").Select() | Out-Null
